# Update "想去人数" (F column) values for matching rows on the
# "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# Map of row -> new F-column value, identical on both affected sheets.
$updates = @{
    6  = 518
    11 = 4148
    13 = 271
    18 = 2858
    29 = 187
    30 = 285
    31 = 1627
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
